$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5532988141424
$ws.Range("B3").Value = 5532988141424
$ws.Range("B3").Select()
